$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Colour the "All Tasks page" detail line (FFC000 / gold) starting at
#    "- Show maybe 3-4 random ..." through to the end of that paragraph.
# ---------------------------------------------------------------------
$search = $d.Range(0, $d.Content.End)
$found = $search.Find.Execute("- Show maybe 3-4 random task details")
if ($found) {
    $startPos = $search.Start

    $paraRange = $d.Range($startPos, $startPos)
    [void]$paraRange.Expand(4)      # wdParagraph -- grab the whole paragraph
    $endPos = $paraRange.End - 1    # exclude the paragraph mark

    $colourRange = $d.Range($startPos, $endPos)
    $colourRange.Font.Color = 49407   # wdColor value for RGB FFC000
}

# ---------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark off the "User's tasks" paragraph and
#    onto the (now-last) empty paragraph at the end of the document --
#    this is what Word leaves behind after the most recent edit.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last

# The very last paragraph has no run/content node yet, which this COM
# surface needs in order to anchor a new bookmark -- so seed it with a
# placeholder character, plant "_GoBack" there (this automatically
# relocates the single allowed "_GoBack" bookmark), then remove the
# placeholder again, leaving just the empty bookmarked paragraph.
$lastPara.Range.InsertBefore("x")
$anchor = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$anchor.Bookmarks.Add("_GoBack")
$placeholder = $d.Range($lastPara.Range.Start, $lastPara.Range.Start + 1)
$placeholder.Delete()
